# Generate Report for handoff
#
# Marks the "879e6e22-dc92-4e1c-823e-b54891410cd1.md" file row as ready for
# handoff across the Overview, zh-cn and de-de sheets, and stamps the new
# handoff datetimes for each locale.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-18 02:51:40"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-18 02:51:50"
